{"js": "// Replace the worksheet date and all the three-digit x one-digit\n// multiplication prompts with the new values from the commit.\nconst replacements = [\n  [\"2024-07-19 Friday\", \"2024-07-20 Saturday\"],\n  [\"230\u00d73=\", \"169\u00d77=\"],\n  [\"641\u00d73=\", \"196\u00d78=\"],\n  [\"308\u00d79=\", \"141\u00d73=\"],\n  [\"291\u00d79=\", \"659\u00d76=\"],\n  [\"510\u00d72=\", \"850\u00d75=\"],\n  [\"364\u00d76=\", \"804\u00d79=\"],\n  [\"131\u00d72=\", \"617\u00d72=\"],\n  [\"231\u00d79=\", \"962\u00d75=\"],\n  [\"542\u00d72=\", \"555\u00d79=\"],\n  [\"605\u00d79=\", \"331\u00d75=\"],\n  [\"127\u00d76=\", \"831\u00d73=\"],\n  [\"472\u00d75=\", \"712\u00d74=\"],\n  [\"770\u00d79=\", \"461\u00d73=\"],\n  [\"413\u00d75=\", \"705\u00d78=\"],\n  [\"693\u00d75=\", \"439\u00d78=\"],\n  [\"593\u00d78=\", \"270\u00d76=\"],\n  [\"495\u00d72=\", \"828\u00d79=\"],\n  [\"909\u00d77=\", \"253\u00d78=\"],\n  [\"767\u00d75=\", \"740\u00d79=\"],\n  [\"221\u00d74=\", \"844\u00d72=\"],\n  [\"942\u00d72=\", \"706\u00d77=\"],\n  [\"857\u00d73=\", \"836\u00d75=\"],\n  [\"461\u00d75=\", \"556\u00d74=\"],\n  [\"113\u00d78=\", \"513\u00d77=\"],\n  [\"870\u00d79=\", \"736\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and all the three-digit x one-digit\n# multiplication prompts with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-19 Friday\", \"2024-07-20 Saturday\"),\n    @(\"230\u00d73=\", \"169\u00d77=\"),\n    @(\"641\u00d73=\", \"196\u00d78=\"),\n    @(\"308\u00d79=\", \"141\u00d73=\"),\n    @(\"291\u00d79=\", \"659\u00d76=\"),\n    @(\"510\u00d72=\", \"850\u00d75=\"),\n    @(\"364\u00d76=\", \"804\u00d79=\"),\n    @(\"131\u00d72=\", \"617\u00d72=\"),\n    @(\"231\u00d79=\", \"962\u00d75=\"),\n    @(\"542\u00d72=\", \"555\u00d79=\"),\n    @(\"605\u00d79=\", \"331\u00d75=\"),\n    @(\"127\u00d76=\", \"831\u00d73=\"),\n    @(\"472\u00d75=\", \"712\u00d74=\"),\n    @(\"770\u00d79=\", \"461\u00d73=\"),\n    @(\"413\u00d75=\", \"705\u00d78=\"),\n    @(\"693\u00d75=\", \"439\u00d78=\"),\n    @(\"593\u00d78=\", \"270\u00d76=\"),\n    @(\"495\u00d72=\", \"828\u00d79=\"),\n    @(\"909\u00d77=\", \"253\u00d78=\"),\n    @(\"767\u00d75=\", \"740\u00d79=\"),\n    @(\"221\u00d74=\", \"844\u00d72=\"),\n    @(\"942\u00d72=\", \"706\u00d77=\"),\n    @(\"857\u00d73=\", \"836\u00d75=\"),\n    @(\"461\u00d75=\", \"556\u00d74=\"),\n    @(\"113\u00d78=\", \"513\u00d77=\"),\n    @(\"870\u00d79=\", \"736\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
